$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.019.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.09%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.72"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.71%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -8.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0962"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.11"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.35"
$ws.Range("D12").ClearFormats()

$ws.Range("E13").Value = "  -2.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -8.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.561.05"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.850"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.228.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.889.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -8.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.53%  "

$ws.Range("E24").Value = "  +12.23%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.44"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.83"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.30"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -8.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.15%  "

$ws.Range("E33").Value = "  -5.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0709"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.65"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +23.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0279"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "66.55"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.81"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -12.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0999"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.189"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.09%  "

$ws.Range("E49").Value = "  -3.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.53%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.10"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.46%  "
